# Insert two new columns (sum_SASA, max_SASA) between the existing
# "SASA" column (C) and "flexibility" column (D), shifting everything
# from D onward two columns to the right (D:E become F:G, etc.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("D:E").Insert()

# New header labels
$ws.Range("D1").Value = "sum_SASA"
$ws.Range("E1").Value = "max_SASA"

# New column values (sum_SASA / max_SASA) per row
$ws.Range("D2").Value = 5.6182815656425
$ws.Range("E2").Value = 3.041646710200382

$ws.Range("D3").Value = 5.318863841181464
$ws.Range("E3").Value = 2.673333100810235

$ws.Range("D4").Value = 5.192981897526153
$ws.Range("E4").Value = 2.597460735209333

$ws.Range("D5").Value = 5.090388917030758
$ws.Range("E5").Value = 2.595288970498892

$ws.Range("D6").Value = 4.455162269858452
$ws.Range("E6").Value = 2.528652387176547

$ws.Range("D7").Value = 4.60605092561061
$ws.Range("E7").Value = 2.340246256923699

$ws.Range("D8").Value = 4.275544050447586
$ws.Range("E8").Value = 2.254148443376134

$ws.Range("D9").Value = 4.479607505862927
$ws.Range("E9").Value = 2.302631596004201

$ws.Range("D10").Value = 4.541613194916014
$ws.Range("E10").Value = 2.305444953997026

$ws.Range("D11").Value = 4.780949935970851
$ws.Range("E11").Value = 2.467945786510408

$ws.Range("D12").Value = 4.603195718807239
$ws.Range("E12").Value = 2.368503450093663

$ws.Range("D13").Value = 4.822099857651121
$ws.Range("E13").Value = 2.6094836119382
